# Bugfix for the naive forecaster component module:
# column A held raw Excel date serials (quarterly anchor dates) formatted
# with a custom "YYYY-MM-DD HH:MM:SS" number format. Replace those with
# plain "<year>Q<quarter>" text labels (e.g. 2005Q1) sharing the same
# bordered/centered header style as A1/B1, and drop the now-unused
# date number format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Figure out the used range on column A (header in row 1, data below).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Copy the header's style (border + bold + centered, no special number
# format) onto every data cell in column A before we overwrite the
# values, so the resulting cells pick up style index 1 (same as A1/B1)
# instead of keeping the old date-formatted style.
$ws.Range("A1").Copy()
$ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item($lastRow, 1)).PasteSpecial(-4122)
$excel.CutCopyMode = 0

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $serial = $cell.Value2
    $dt = [DateTime]::FromOADate($serial)
    $quarter = [Math]::Floor(($dt.Month - 1) / 3) + 1
    $cell.Value2 = "$($dt.Year)Q$quarter"
}
